$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Chole chawal" row (old row 3) so it becomes row 19; rows that
# used to be 4-19 shift up by one to become rows 3-18.
$ws.Rows("3:3").Delete()
$ws.Rows("19:19").Insert()

$ws.Range("A19").Value2 = "Chole chawal"
$ws.Range("B19").Value2 = 50
$ws.Range("C19").Value2 = 80
$ws.Range("D19").Value2 = "Chole Chawal.jpeg"

# Match the resulting selection/view state from the edit.
$ws.Range("D20").Select()
